$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 7.065029
$ws.Range("H2").Value2 = 21.195087
$ws.Range("I2").Value2 = 0.2168972219923977
$ws.Range("J2").Value2 = 0.2168972219923977
$ws.Range("M2").Value2 = 0.5001966666666666
$ws.Range("N2").Value2 = 1.50059
$ws.Range("O2").Value2 = 0.03894027965151046
$ws.Range("P2").Value2 = 0.03894027965151046
$ws.Range("Q2").Value2 = 3.533903955703333
$ws.Range("R2").Value2 = 31.80513560133
$ws.Range("S2").Value2 = 0.008446038480019714
$ws.Range("T2").Value2 = 0.008446038480019713

$ws.Range("G3").Value2 = 7.065029
$ws.Range("H3").Value2 = 21.195087
$ws.Range("I3").Value2 = 0.2168972219923977
$ws.Range("J3").Value2 = 0.2168972219923977
$ws.Range("O3").Value2 = 0.7732779360092192
$ws.Range("P3").Value2 = 0.7732779360092191
$ws.Range("Q3").Value2 = 70.17643379495065
$ws.Range("R3").Value2 = 631.587904154556
$ws.Range("S3").Value2 = 0.1677218361484148
$ws.Range("T3").Value2 = 0.1677218361484147

$ws.Range("G4").Value2 = 7.065029
$ws.Range("H4").Value2 = 21.195087
$ws.Range("I4").Value2 = 0.2168972219923977
$ws.Range("J4").Value2 = 0.2168972219923977
$ws.Range("M4").Value2 = 2.334238666666666
$ws.Range("N4").Value2 = 7.002715999999999
$ws.Range("O4").Value2 = 0.1817203362411497
$ws.Range("P4").Value2 = 0.1817203362411496
$ws.Range("Q4").Value2 = 16.49146387292133
$ws.Range("R4").Value2 = 148.423174856292
$ws.Range("S4").Value2 = 0.0394146361102298
$ws.Range("T4").Value2 = 0.03941463611022979

$ws.Range("G5").Value2 = 7.065029
$ws.Range("H5").Value2 = 21.195087
$ws.Range("I5").Value2 = 0.2168972219923977
$ws.Range("J5").Value2 = 0.2168972219923977
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.07786066666666666
$ws.Range("N5").Value2 = 0.233582
$ws.Range("O5").Value2 = 0.006061448098120818
$ws.Range("P5").Value2 = 0.006061448098120817
$ws.Range("Q5").Value2 = 0.5500878679593333
$ws.Range("R5").Value2 = 4.950790811634
$ws.Range("S5").Value2 = 0.001314711253733508
$ws.Range("T5").Value2 = 0.001314711253733508

$ws.Range("I6").Value2 = 0.4849997229314377
$ws.Range("J6").Value2 = 0.4849997229314376
$ws.Range("M6").Value2 = 0.5001966666666666
$ws.Range("N6").Value2 = 1.50059
$ws.Range("O6").Value2 = 0.03894027965151046
$ws.Range("P6").Value2 = 0.03894027965151046
$ws.Range("Q6").Value2 = 7.902094935279999
$ws.Range("R6").Value2 = 71.11885441752
$ws.Range("S6").Value2 = 0.01888602484185527
$ws.Range("T6").Value2 = 0.01888602484185527

$ws.Range("I7").Value2 = 0.4849997229314377
$ws.Range("J7").Value2 = 0.4849997229314376
$ws.Range("O7").Value2 = 0.7732779360092192
$ws.Range("P7").Value2 = 0.7732779360092191
$ws.Range("S7").Value2 = 0.3750395847134653
$ws.Range("T7").Value2 = 0.3750395847134652

$ws.Range("I8").Value2 = 0.4849997229314377
$ws.Range("J8").Value2 = 0.4849997229314376
$ws.Range("M8").Value2 = 2.334238666666666
$ws.Range("N8").Value2 = 7.002715999999999
$ws.Range("O8").Value2 = 0.1817203362411497
$ws.Range("P8").Value2 = 0.1817203362411496
$ws.Range("Q8").Value2 = 36.876246434272
$ws.Range("R8").Value2 = 331.886217908448
$ws.Range("S8").Value2 = 0.08813431272796528
$ws.Range("T8").Value2 = 0.08813431272796525

$ws.Range("I9").Value2 = 0.4849997229314377
$ws.Range("J9").Value2 = 0.4849997229314376
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 0.07786066666666666
$ws.Range("N9").Value2 = 0.233582
$ws.Range("O9").Value2 = 0.006061448098120818
$ws.Range("P9").Value2 = 0.006061448098120817
$ws.Range("Q9").Value2 = 1.230040943344
$ws.Range("R9").Value2 = 11.070368490096
$ws.Range("S9").Value2 = 0.002939800648151887
$ws.Range("T9").Value2 = 0.002939800648151885

$ws.Range("G10").Value2 = 4.745838333333333
$ws.Range("H10").Value2 = 14.237515
$ws.Range("I10").Value2 = 0.1456977955115303
$ws.Range("J10").Value2 = 0.1456977955115302
$ws.Range("M10").Value2 = 0.5001966666666666
$ws.Range("N10").Value2 = 1.50059
$ws.Range("O10").Value2 = 0.03894027965151046
$ws.Range("P10").Value2 = 0.03894027965151046
$ws.Range("Q10").Value2 = 2.373852514872222
$ws.Range("R10").Value2 = 21.36467263385
$ws.Range("S10").Value2 = 0.005673512901827574
$ws.Range("T10").Value2 = 0.005673512901827572

$ws.Range("G11").Value2 = 4.745838333333333
$ws.Range("H11").Value2 = 14.237515
$ws.Range("I11").Value2 = 0.1456977955115303
$ws.Range("J11").Value2 = 0.1456977955115302
$ws.Range("O11").Value2 = 0.7732779360092192
$ws.Range("P11").Value2 = 0.7732779360092191
$ws.Range("Q11").Value2 = 47.14007679242444
$ws.Range("R11").Value2 = 424.26069113182
$ws.Range("S11").Value2 = 0.1126648905942494
$ws.Range("T11").Value2 = 0.1126648905942493

$ws.Range("G12").Value2 = 4.745838333333333
$ws.Range("H12").Value2 = 14.237515
$ws.Range("I12").Value2 = 0.1456977955115303
$ws.Range("J12").Value2 = 0.1456977955115302
$ws.Range("M12").Value2 = 2.334238666666666
$ws.Range("N12").Value2 = 7.002715999999999
$ws.Range("O12").Value2 = 0.1817203362411497
$ws.Range("P12").Value2 = 0.1817203362411496
$ws.Range("Q12").Value2 = 11.07791934341555
$ws.Range("R12").Value2 = 99.70127409074
$ws.Range("S12").Value2 = 0.02647625238994954
$ws.Range("T12").Value2 = 0.02647625238994953

$ws.Range("G13").Value2 = 4.745838333333333
$ws.Range("H13").Value2 = 14.237515
$ws.Range("I13").Value2 = 0.1456977955115303
$ws.Range("J13").Value2 = 0.1456977955115302
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 0.07786066666666666
$ws.Range("N13").Value2 = 0.233582
$ws.Range("O13").Value2 = 0.006061448098120818
$ws.Range("P13").Value2 = 0.006061448098120817
$ws.Range("Q13").Value2 = 0.3695141365255555
$ws.Range("R13").Value2 = 3.32562722873
$ws.Range("S13").Value2 = 0.0008831396255037609
$ws.Range("T13").Value2 = 0.0008831396255037603

$ws.Range("G14").Value2 = 4.964321666666666
$ws.Range("H14").Value2 = 14.892965
$ws.Range("I14").Value2 = 0.1524052595646345
$ws.Range("J14").Value2 = 0.1524052595646344
$ws.Range("M14").Value2 = 0.5001966666666666
$ws.Range("N14").Value2 = 1.50059
$ws.Range("O14").Value2 = 0.03894027965151046
$ws.Range("P14").Value2 = 0.03894027965151046
$ws.Range("Q14").Value2 = 2.483137149927777
$ws.Range("R14").Value2 = 22.34823434934999
$ws.Range("S14").Value2 = 0.005934703427807906
$ws.Range("T14").Value2 = 0.005934703427807905

$ws.Range("G15").Value2 = 4.964321666666666
$ws.Range("H15").Value2 = 14.892965
$ws.Range("I15").Value2 = 0.1524052595646345
$ws.Range("J15").Value2 = 0.1524052595646344
$ws.Range("O15").Value2 = 0.7732779360092192
$ws.Range("P15").Value2 = 0.7732779360092191
$ws.Range("Q15").Value2 = 49.31025630293554
$ws.Range("R15").Value2 = 443.79230672642
$ws.Range("S15").Value2 = 0.1178516245530899
$ws.Range("T15").Value2 = 0.1178516245530898

$ws.Range("G16").Value2 = 4.964321666666666
$ws.Range("H16").Value2 = 14.892965
$ws.Range("I16").Value2 = 0.1524052595646345
$ws.Range("J16").Value2 = 0.1524052595646344
$ws.Range("M16").Value2 = 2.334238666666666
$ws.Range("N16").Value2 = 7.002715999999999
$ws.Range("O16").Value2 = 0.1817203362411497
$ws.Range("P16").Value2 = 0.1817203362411496
$ws.Range("Q16").Value2 = 11.58791158810444
$ws.Range("R16").Value2 = 104.29120429294
$ws.Range("S16").Value2 = 0.02769513501300506
$ws.Range("T16").Value2 = 0.02769513501300506

$ws.Range("G17").Value2 = 4.964321666666666
$ws.Range("H17").Value2 = 14.892965
$ws.Range("I17").Value2 = 0.1524052595646345
$ws.Range("J17").Value2 = 0.1524052595646344
$ws.Range("K17").Value2 = 3
$ws.Range("L17").Value2 = 1
$ws.Range("M17").Value2 = 0.07786066666666666
$ws.Range("N17").Value2 = 0.233582
$ws.Range("O17").Value2 = 0.006061448098120818
$ws.Range("P17").Value2 = 0.006061448098120817
$ws.Range("Q17").Value2 = 0.3865253945144443
$ws.Range("R17").Value2 = 3.478728550629999
$ws.Range("S17").Value2 = 0.0009237965707316632
$ws.Range("T17").Value2 = 0.0009237965707316632
